$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.008.34"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3
$ws.Range("D3").Value = "1.584.18"
$ws.Range("E3").Value = "  -1.87%  "

# Row 4
$ws.Range("E4").Value = "  +0.49%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "300.45"
$ws.Range("E6").Value = "  -0.83%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3753"
$ws.Range("E7").Value = "  -0.80%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3569"
$ws.Range("E8").Value = "  -3.29%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.46"
$ws.Range("E9").Value = "  +2.72%  "

# Row 10
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.212"
$ws.Range("E11").Value = "  -5.68%  "

# Row 12
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07970"
$ws.Range("E12").Value = "  -1.82%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.74"
$ws.Range("E13").Value = "  -6.58%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.437"
$ws.Range("E14").Value = "  -3.41%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.266"
$ws.Range("E15").Value = "  -5.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001217"
$ws.Range("E16").Value = "  -4.84%  "

# Row 17
$ws.Range("D17").Value = "1.587.03"
$ws.Range("E17").Value = "  -1.37%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.09"
$ws.Range("E18").Value = "  +0.33%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06756"
$ws.Range("E19").Value = "  -0.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  -4.31%  "

# Row 21
$ws.Range("E21").Value = "  +0.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.387"
$ws.Range("E22").Value = "  -3.53%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.71"
$ws.Range("E23").Value = "  -3.09%  "

# Row 24
$ws.Range("D24").Value = "23.007.56"
$ws.Range("E24").Value = "  -0.90%  "

# Row 25
$ws.Range("E25").Value = "  +0.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.747"
$ws.Range("E26").Value = "  -6.58%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.62"
$ws.Range("E27").Value = "  -2.69%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "146.62"
$ws.Range("E28").Value = "  -2.90%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.197"
$ws.Range("E29").Value = "  -1.71%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.73"
$ws.Range("E30").Value = "  -0.80%  "

# Row 31
$ws.Range("E31").Value = "  -4.93%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.488"
$ws.Range("E32").Value = "  -7.95%  "

# Row 33
$ws.Range("D33").Value = "1.765.06"
$ws.Range("E33").Value = "  -1.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9289"
$ws.Range("E34").Value = "  -7.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07303"
$ws.Range("E35").Value = "  -6.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02657"
$ws.Range("E36").Value = "  -5.22%  "

# Row 37
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08729"
$ws.Range("E37").Value = "  -1.88%  "

# Row 38
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2477"
$ws.Range("E38").Value = "  -3.59%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.876"
$ws.Range("E39").Value = "  -3.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.992"
$ws.Range("E40").Value = "  -5.83%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.330"
$ws.Range("E41").Value = "  -5.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6819"
$ws.Range("E42").Value = "  -5.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.75"
$ws.Range("E43").Value = "  -8.81%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.67"
$ws.Range("E44").Value = "  -8.52%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6334"
$ws.Range("E45").Value = "  -5.03%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.976"
$ws.Range("E46").Value = "  -0.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.234"
$ws.Range("E47").Value = "  -3.78%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.81"
$ws.Range("E48").Value = "  -0.78%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07873"
$ws.Range("E49").Value = "  -2.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.180"
$ws.Range("E50").Value = "  +0.26%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.174"
$ws.Range("E51").Value = "  -1.39%  "
